$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 34, pushing the current row 34
# (Granada / Segunda / 2022-04-13) down to row 36.
$ws.Range("A34:A35").EntireRow.Insert()

# The new rows 34 and 35 keep the "old week" (2022-04-13) data that used
# to live in rows 32 and 33 (Especial / Primera, Provincia de Limari).
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34, 3).Value = "Metropolitana"
$ws.Cells.Item(34, 4).Value = 44664
$ws.Cells.Item(34, 5).Value = 13
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100104
$ws.Cells.Item(34, 8).Value = "Frutos de pepita"
$ws.Cells.Item(34, 9).Value = 100104001
$ws.Cells.Item(34, 10).Value = "Granada"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Especial"
$ws.Cells.Item(34, 13).Value = 300
$ws.Cells.Item(34, 14).Value = 21600
$ws.Cells.Item(34, 15).Value = 21600
$ws.Cells.Item(34, 16).Value = 21600
$ws.Cells.Item(34, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(34, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 19).Value = 1200
$ws.Cells.Item(34, 20).Value = 18

$ws.Cells.Item(35, 1).Value = 9
$ws.Cells.Item(35, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(35, 3).Value = "Metropolitana"
$ws.Cells.Item(35, 4).Value = 44664
$ws.Cells.Item(35, 5).Value = 13
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100104
$ws.Cells.Item(35, 8).Value = "Frutos de pepita"
$ws.Cells.Item(35, 9).Value = 100104001
$ws.Cells.Item(35, 10).Value = "Granada"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 250
$ws.Cells.Item(35, 14).Value = 18000
$ws.Cells.Item(35, 15).Value = 18000
$ws.Cells.Item(35, 16).Value = 18000
$ws.Cells.Item(35, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(35, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(35, 19).Value = 1000
$ws.Cells.Item(35, 20).Value = 18

# Rows 32 and 33 now hold the new week's prices (2022-05-25, Wonderfull
# variety, Region de O'Higgins).
$ws.Cells.Item(32, 4).Value = 44706
$ws.Cells.Item(32, 11).Value = "Wonderfull"
$ws.Cells.Item(32, 13).Value = 200
$ws.Cells.Item(32, 14).Value = 16000
$ws.Cells.Item(32, 15).Value = 16000
$ws.Cells.Item(32, 16).Value = 16000
$ws.Cells.Item(32, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 19).Value = 889

$ws.Cells.Item(33, 4).Value = 44706
$ws.Cells.Item(33, 11).Value = "Wonderfull"
$ws.Cells.Item(33, 13).Value = 220
$ws.Cells.Item(33, 14).Value = 12500
$ws.Cells.Item(33, 15).Value = 12500
$ws.Cells.Item(33, 16).Value = 12500
$ws.Cells.Item(33, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(33, 19).Value = 694
